$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 41558.4
$ws.Range("J3").Value = 41558.4
$ws.Range("L3").Value = 41558.4
$ws.Range("N3").Value = -41786.4
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H58").Value = 1343.3334
$ws.Range("J58").Value = 1800
$ws.Range("L58").Value = 5400
$ws.Range("N58").Value = -5700
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H102").Value = 41558.4
$ws.Range("J102").Value = 41558.4
$ws.Range("L102").Value = 41558.4
$ws.Range("N102").Value = -48048.4
$ws.Range("H113").Value = 2999
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 6355.913
$ws.Range("I116").Value = 4047.5715
$ws.Range("J116").Value = 7365.8125
$ws.Range("K116").Value = 4047.5715
$ws.Range("L116").Value = 7365.8125
$ws.Range("M116").Value = -605.5715
$ws.Range("N116").Value = -14249.8125
$ws.Range("H132").Value = 18260.09
$ws.Range("I132").Value = 2908.5881
$ws.Range("K132").Value = 8725.764299999999
$ws.Range("M132").Value = -6195.764299999999
$ws.Range("H137").Value = 3099.7214
$ws.Range("I137").Value = 1250.3572
$ws.Range("J137").Value = 3498.0461
$ws.Range("K137").Value = 3751.0716
$ws.Range("L137").Value = 10494.1383
$ws.Range("M137").Value = -1201.0716
$ws.Range("N137").Value = -15594.1383
$ws.Range("H138").Value = 2509.2341
$ws.Range("I138").Value = 1730.721
$ws.Range("J138").Value = 3165.6274
$ws.Range("K138").Value = 5192.163
$ws.Range("L138").Value = 9496.8822
$ws.Range("M138").Value = -52.16300000000047
$ws.Range("N138").Value = -19776.8822
$ws.Range("H141").Value = 2558.0278
$ws.Range("I141").Value = 1559.1333
$ws.Range("J141").Value = 7552.5
$ws.Range("K141").Value = 4677.3999
$ws.Range("L141").Value = 22657.5
$ws.Range("M141").Value = 502.6000999999997
$ws.Range("N141").Value = -33017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 385422.3
$ws.Range("I5").Value = 417349.16
$ws.Range("J5").Value = 2300
$ws.Range("K5").Value = 417349.16
$ws.Range("L5").Value = 2300
$ws.Range("M5").Value = -417237.16
$ws.Range("N5").Value = -2524
$ws.Range("H32").Value = 20630.06
$ws.Range("I32").Value = 19828.344
$ws.Range("J32").Value = 100000
$ws.Range("K32").Value = 19828.344
$ws.Range("L32").Value = 100000
$ws.Range("M32").Value = -19541.344
$ws.Range("N32").Value = -100574
$ws.Range("H61").Value = 1528.1765
$ws.Range("I61").Value = 1332.8478
$ws.Range("J61").Value = 1936.591
$ws.Range("K61").Value = 1332.8478
$ws.Range("L61").Value = 1936.591
$ws.Range("M61").Value = -1120.8478
$ws.Range("N61").Value = -2360.591
$ws.Range("H114").Value = 33775.5
$ws.Range("J114").Value = 33775.5
$ws.Range("L114").Value = 33775.5
$ws.Range("N114").Value = -42453.5
$ws.Range("H132").Value = 7248434.5
$ws.Range("I132").Value = 12196437
$ws.Range("J132").Value = 3144.4285
$ws.Range("K132").Value = 36589311
$ws.Range("L132").Value = 9433.2855
$ws.Range("M132").Value = -36586781
$ws.Range("N132").Value = -14493.2855
$ws.Range("H133").Value = 29987.615
$ws.Range("J133").Value = 29987.615
$ws.Range("L133").Value = 29987.615
$ws.Range("N133").Value = -35047.61500000001
$ws.Range("H136").Value = 1528.1765
$ws.Range("I136").Value = 1332.8478
$ws.Range("J136").Value = 1936.591
$ws.Range("K136").Value = 3998.5434
$ws.Range("L136").Value = 5809.772999999999
$ws.Range("M136").Value = -1448.5434
$ws.Range("N136").Value = -10909.773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 385422.3
$ws.Range("I4").Value = 417349.16
$ws.Range("J4").Value = 2300
$ws.Range("K4").Value = 417349.16
$ws.Range("L4").Value = 2300
$ws.Range("M4").Value = -417234.16
$ws.Range("N4").Value = -2530
$ws.Range("H134").Value = 2130.375
$ws.Range("I134").Value = 1003.8125
$ws.Range("K134").Value = 3011.4375
$ws.Range("M134").Value = -476.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7496.3438
$ws.Range("I31").Value = 3471
$ws.Range("J31").Value = 9326.046
$ws.Range("K31").Value = 3471
$ws.Range("L31").Value = 9326.046
$ws.Range("M31").Value = -3176
$ws.Range("N31").Value = -9916.046
$ws.Range("H34").Value = 7496.3438
$ws.Range("I34").Value = 3471
$ws.Range("J34").Value = 9326.046
$ws.Range("K34").Value = 3471
$ws.Range("L34").Value = 9326.046
$ws.Range("M34").Value = -3269
$ws.Range("N34").Value = -9730.046
$ws.Range("H99").Value = 2192.923
$ws.Range("I99").Value = 2300.7273
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 2300.7273
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -802.7273
$ws.Range("N99").Value = -4596
$ws.Range("H126").Value = 2192.923
$ws.Range("I126").Value = 2300.7273
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 6902.1819
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -4432.1819
$ws.Range("N126").Value = -9740
$ws.Range("H134").Value = 2202.7144
$ws.Range("I134").Value = 1350.4445
$ws.Range("J134").Value = 7316.3335
$ws.Range("K134").Value = 4051.3335
$ws.Range("L134").Value = 21949.0005
$ws.Range("M134").Value = -1516.3335
$ws.Range("N134").Value = -27019.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9563853
$ws.Range("I2").Value = 59
$ws.Range("J2").Value = 16138962
$ws.Range("K2").Value = 354
$ws.Range("L2").Value = 96833772
$ws.Range("M2").Value = -241
$ws.Range("N2").Value = -96833998
$ws.Range("H4").Value = 921.5454999999999
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 3112.3333
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 9336.999899999999
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -9560.999899999999
$ws.Range("H24").Value = 1497.5
$ws.Range("I24").Value = 290
$ws.Range("J24").Value = 1900
$ws.Range("K24").Value = 870
$ws.Range("L24").Value = 5700
$ws.Range("M24").Value = -640
$ws.Range("N24").Value = -6160
$ws.Range("H133").Value = 7279.077
$ws.Range("I133").Value = 9246
$ws.Range("J133").Value = 6049.75
$ws.Range("K133").Value = 27738
$ws.Range("L133").Value = 18149.25
$ws.Range("M133").Value = -22678
$ws.Range("N133").Value = -28269.25
$ws.Range("H140").Value = 1761.6571
$ws.Range("I140").Value = 1568.1428
$ws.Range("J140").Value = 2535.7144
$ws.Range("K140").Value = 4704.428400000001
$ws.Range("L140").Value = 7607.1432
$ws.Range("M140").Value = 475.5715999999993
$ws.Range("N140").Value = -17967.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1854.6666
$ws.Range("I102").Value = 550
$ws.Range("J102").Value = 2507
$ws.Range("K102").Value = 550
$ws.Range("L102").Value = 2507
$ws.Range("M102").Value = 1072
$ws.Range("N102").Value = -5751
$ws.Range("H126").Value = 15346.934
$ws.Range("I126").Value = 21780.4
$ws.Range("K126").Value = 65341.2
$ws.Range("M126").Value = -62871.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 399.5
$ws.Range("I22").Value = 399.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 399.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -104.5
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 399.5
$ws.Range("I27").Value = 399.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 399.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -292.5
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 2593.5
$ws.Range("I40").Value = 2499.7856
$ws.Range("J40").Value = 3249.5
$ws.Range("K40").Value = 2499.7856
$ws.Range("L40").Value = 3249.5
$ws.Range("M40").Value = -2363.7856
$ws.Range("N40").Value = -3521.5
$ws.Range("H94").Value = 60397.332
$ws.Range("J94").Value = 60397.332
$ws.Range("L94").Value = 60397.332
$ws.Range("N94").Value = -61749.332
$ws.Range("H122").Value = 28887.676
$ws.Range("I122").Value = 29651.084
$ws.Range("J122").Value = 1405
$ws.Range("K122").Value = 88953.25199999999
$ws.Range("L122").Value = 4215
$ws.Range("M122").Value = -86503.25199999999
$ws.Range("N122").Value = -9115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1166.6666
$ws.Range("J2").Value = 1166.6666
$ws.Range("L2").Value = 1166.6666
$ws.Range("N2").Value = -1390.6666
$ws.Range("H122").Value = 35238876
$ws.Range("I122").Value = 44048344
$ws.Range("J122").Value = 998.3333
$ws.Range("K122").Value = 132145032
$ws.Range("L122").Value = 2994.9999
$ws.Range("M122").Value = -132142582
$ws.Range("N122").Value = -7894.9999
$ws.Range("H126").Value = 2674615.8
$ws.Range("I126").Value = 5882754.5
$ws.Range("K126").Value = 17648263.5
$ws.Range("M126").Value = -17645793.5
$ws.Range("H136").Value = 17634.143
$ws.Range("I136").Value = 33443.902
$ws.Range("K136").Value = 100331.706
$ws.Range("M136").Value = -97781.70600000001
